$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.08856766666666667
$ws.Range("H2").Value = 0.265703
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 10.41202966666667
$ws.Range("N2").Value = 31.236089
$ws.Range("Q2").Value = 0.9221691728407778
$ws.Range("R2").Value = 8.299522555567
